$d = $word.ActiveDocument

# Merge the fragmented runs in the Title paragraph into a single run.
$d.Content.Find.Execute("Questions: Introduction to radians", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Questions: Introduction to radians", 2)

# Merge the fragmented runs in the Author paragraph into a single run.
$d.Content.Find.Execute("Mark Toner, Ifan Howell-Baines", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mark Toner, Ifan Howell-Baines", 2)

# Merge the fragmented runs in the Abstract paragraph into a single run.
$d.Content.Find.Execute("Questions relating to the introduction to radians study guide.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Questions relating to the introduction to radians study guide.", 2)
